# Chapter2-1a.xlsx edit:
# - Clear the leftover "SUM(X)/SUM(Y)/SUM(XY)/SUM(X^2)/n/m/b" labels that were
#   sitting in C3:C9 (their styling stays, only the text content is removed).
#   This shrinks the shared-string table and the surviving strings get
#   renumbered automatically by Excel - no further action needed for that.
# - Move the active selection on the sheet from H12 to H19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray formula-name labels in column C (rows 3 through 9).
$ws.Range("C3:C9").ClearContents()

# Update the selected cell shown in the saved sheet view.
$ws.Range("H19").Select()
